# Update "想去人数" (people-who-want-to-go) counts across the workbook.
# Mapping of sheets: 展览=1, 演出=2, 本地生活=3, 全部类型=4

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 5743
$wsExpo.Range("F4").Value = 93
$wsExpo.Range("F5").Value = 408

# Sheet "演出"
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 63

# Sheet "全部类型" (combined listing)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 5743
$wsAll.Range("F4").Value = 93
$wsAll.Range("F5").Value = 63
$wsAll.Range("F6").Value = 408
